$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New integer reward values (rows 2-17, columns B-E), replacing the old decimal values
$values = @{
    2  = @(0, 1, 10, 0)
    3  = @(0, 5, 20, 0)
    4  = @(0, 0, 0, 20)
    5  = @(0, 0, 0, 0)
    6  = @(0, 0, 1, 0)
    7  = @(0, 0, 20, 0)
    8  = @(0, 0, 0, 20)
    9  = @(0, 0, 0, 0)
    10 = @(0, 1, 1, 0)
    11 = @(0, 5, 20, 0)
    12 = @(0, 0, 0, 20)
    13 = @(0, 0, 0, 0)
    14 = @(0, 1, 10, 0)
    15 = @(0, 5, 20, 0)
    16 = @(0, 0, 0, 20)
    17 = @(0, 0, 0, 0)
}

foreach ($row in $values.Keys) {
    $rowVals = $values[$row]
    for ($col = 2; $col -le 5; $col++) {
        $cell = $ws.Cells.Item($row, $col)
        $cell.Value = $rowVals[$col - 2]
        $cell.NumberFormat = "0"
    }
}

# Update the selected cell to match the final saved state
$ws.Range("J14").Select()
